# Applies the cryptos price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Prefix with an apostrophe so Excel stores the content as literal text
    # (prevents auto-conversion of numeric-looking strings into numbers),
    # then reset the style so no extra quote-prefix/number-format is left behind.
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextCell "D2" "68.321.27"
Set-TextCell "E2" "  +1.58%  "

Set-TextCell "D3" "2.639.51"
Set-TextCell "E3" "  +1.24%  "

Set-TextCell "E4" "  +0.02%  "

Set-TextCell "D5" "599.15"
Set-TextCell "E5" "  +1.29%  "

Set-TextCell "D6" "154.05"
Set-TextCell "E6" "  +2.33%  "

Set-TextCell "E7" "  +0.01%  "

Set-TextCell "D8" "0.546"
Set-TextCell "E8" "  -0.07%  "

Set-TextCell "D9" "2.638.07"
Set-TextCell "E9" "  +1.30%  "

Set-TextCell "E10" "  +8.31%  "

Set-TextCell "E11" "  -0.45%  "

Set-TextCell "E12" "  +0.93%  "

Set-TextCell "E13" "  +1.25%  "

Set-TextCell "D14" "27.89"
Set-TextCell "E14" "  +2.27%  "

Set-TextCell "E15" "  +3.33%  "

Set-TextCell "D16" "3.128.38"
Set-TextCell "E16" "  +1.57%  "

Set-TextCell "D17" "68.173.21"
Set-TextCell "E17" "  +1.61%  "

Set-TextCell "D18" "2.641.99"
Set-TextCell "E18" "  +1.43%  "

Set-TextCell "D19" "11.42"
Set-TextCell "E19" "  +3.50%  "

Set-TextCell "D20" "364.75"
Set-TextCell "E20" "  -2.34%  "

Set-TextCell "D21" "7.39"
Set-TextCell "E21" "  +0.10%  "

Set-TextCell "D22" "4.26"
Set-TextCell "E22" "  -0.68%  "

Set-TextCell "E23" "  -0.02%  "

Set-TextCell "E24" "  +2.46%  "

Set-TextCell "D25" "73.37"
Set-TextCell "E25" "  -0.12%  "

Set-TextCell "D26" "0.998"
Set-TextCell "E26" "  -0.27%  "

Set-TextCell "D27" "9.95"
Set-TextCell "E27" "  -0.14%  "

Set-TextCell "D28" "2.773.36"
Set-TextCell "E28" "  +1.29%  "

Set-TextCell "E29" "  +5.25%  "

Set-TextCell "E30" "  -0.20%  "

Set-TextCell "D31" "573.06"
Set-TextCell "E31" "  -1.54%  "

Set-TextCell "D32" "1.42"
Set-TextCell "E32" "  +4.17%  "

Set-TextCell "D33" "7.96"
Set-TextCell "E33" "  +3.76%  "

Set-TextCell "E34" "  +2.45%  "

Set-TextCell "E35" "  +2.77%  "

Set-TextCell "D36" "0.999"
Set-TextCell "E36" "  -0.01%  "

Set-TextCell "E37" "  +3.04%  "

Set-TextCell "D38" "160.20"

Set-TextCell "E39" "  +0.85%  "

Set-TextCell "E40" "  +3.69%  "

Set-TextCell "E41" "  +0.80%  "

Set-TextCell "D42" "5.36"
Set-TextCell "E42" "  +2.34%  "

Set-TextCell "E43" "  +3.58%  "

Set-TextCell "E44" "  +2.24%  "

Set-TextCell "D45" "0.0₆0318"
Set-TextCell "E45" "  +11.58%  "

Set-TextCell "E46" "  +0.01%  "

Set-TextCell "D48" "156.94"
Set-TextCell "E48" "  +2.43%  "

Set-TextCell "E49" "  +0.47%  "

Set-TextCell "E50" "  +1.49%  "

Set-TextCell "D51" "21.83"
Set-TextCell "E51" "  +2.23%  "
